$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "88.481.99"
$ws.Range("E2").Value = "  +9.47%  "
$ws.Range("D3").Value = "3.326.56"
$ws.Range("E3").Value = "  +6.26%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'217.31"
$ws.Range("E5").Value = "  +5.82%  "
$ws.Range("D6").Value = "'649.20"
$ws.Range("E6").Value = "  +4.58%  "
$ws.Range("D7").Value = "'0.397"
$ws.Range("E7").Value = "  +41.96%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.601"
$ws.Range("E9").Value = "  +4.36%  "
$ws.Range("D10").Value = "3.320.57"
$ws.Range("E10").Value = "  +5.99%  "
$ws.Range("D11").Value = "'0.584"
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("D12").Value = "'0.0000291"
$ws.Range("E12").Value = "  +16.67%  "
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "'35.44"
$ws.Range("E14").Value = "  +14.22%  "
$ws.Range("D15").Value = "3.936.33"
$ws.Range("E15").Value = "  +6.13%  "
$ws.Range("E16").Value = "  +5.00%  "
$ws.Range("D17").Value = "88.326.63"
$ws.Range("E17").Value = "  +9.08%  "
$ws.Range("D18").Value = "3.329.28"
$ws.Range("E18").Value = "  +6.44%  "
$ws.Range("E19").Value = "  +5.55%  "
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").Value = "'9.60"
$ws.Range("E21").Value = "  +8.16%  "
$ws.Range("D22").Value = "'454.42"
$ws.Range("E22").Value = "  +6.15%  "
$ws.Range("D23").Value = "'5.42"
$ws.Range("E23").Value = "  +7.59%  "
$ws.Range("D24").Value = "'7.44"
$ws.Range("E24").Value = "  +4.30%  "
$ws.Range("D25").Value = "'5.54"
$ws.Range("E25").Value = "  +8.46%  "
$ws.Range("D26").Value = "'12.89"
$ws.Range("E26").Value = "  +19.42%  "
$ws.Range("D27").Value = "3.512.02"
$ws.Range("E27").Value = "  +6.19%  "
$ws.Range("D28").Value = "'78.19"
$ws.Range("E28").Value = "  +3.71%  "
$ws.Range("D29").Value = "'0.213"
$ws.Range("E29").Value = "  +44.87%  "
$ws.Range("D30").Value = "'0.0000133"
$ws.Range("E30").Value = "  +11.15%  "
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "'9.36"
$ws.Range("D33").Value = "'595.47"
$ws.Range("E33").Value = "  +8.07%  "
$ws.Range("E34").Value = "  +9.08%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("E36").Value = "  +7.22%  "
$ws.Range("D37").Value = "'7.18"
$ws.Range("E37").Value = "  +22.49%  "
$ws.Range("D38").Value = "'0.141"
$ws.Range("E38").Value = "  -5.70%  "
$ws.Range("D39").Value = "'23.10"
$ws.Range("E39").Value = "  +2.89%  "
$ws.Range("D40").Value = "'2.15"
$ws.Range("E40").Value = "  +8.81%  "
$ws.Range("D41").Value = "'0.419"
$ws.Range("E41").Value = "  +4.34%  "
$ws.Range("E42").Value = "  +5.55%  "
$ws.Range("D43").Value = "'0.997"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").Value = "'3.13"
$ws.Range("E44").Value = "  +5.25%  "
$ws.Range("D45").Value = "'158.10"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "'1.43"
$ws.Range("E47").Value = "  +9.78%  "
$ws.Range("D48").Value = "'187.45"
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("D49").Value = "'45.95"
$ws.Range("E49").Value = "  +4.61%  "
$ws.Range("D50").Value = "'4.43"
$ws.Range("E50").Value = "  +6.16%  "
$ws.Range("D51").Value = "'0.658"
$ws.Range("E51").Value = "  +6.43%  "
